# Updated symbol list on Sun Dec 25 08:52:06 UTC 2022 with GitHub Actions
#
# This script reproduces the data refresh captured by the diff: a handful of
# "Price" (column D) values were refreshed, and several coin-listing rows
# shifted down by one slot (each row's Coin/Link/Price/Volume data took on
# the values that used to belong to the row above/below it).
#
# Note: column D stores its numbers as *text* (inlineStr) in the workbook,
# e.g. "245.14" rather than the number 245.14. Assigning a numeric-looking
# string straight to Range.Value lets Excel auto-coerce it into a real
# number (losing exact text formatting / precision), so for those cells we
# briefly force the Text number format, assign the string, and then restore
# the cell's original (unstyled) look so no stray style gets left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# A plain, unstyled data cell we can borrow the "Normal" style from once
# we're done forcing Text format, so cells don't end up with a leftover
# explicit style index.
$normalStyle = $ws.Range("G2").Style

function Set-TextCell($range, [string]$value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $normalStyle
}

# ---------------------------------------------------------------------
# Simple "Price" refreshes (column D only)
# ---------------------------------------------------------------------
Set-TextCell $ws.Range("D2")  "245.14"
Set-TextCell $ws.Range("D4")  "5.416"
Set-TextCell $ws.Range("D5")  "0.05993"
Set-TextCell $ws.Range("D6")  "3.392"
Set-TextCell $ws.Range("D7")  "0.8074"
Set-TextCell $ws.Range("D8")  "0.9273"

# ---------------------------------------------------------------------
# Rows 9-17: each row now shows the coin that used to be one row below it
# (the "One" entry moves from row 9 down to row 17), with refreshed prices.
# ---------------------------------------------------------------------

# Row 9 (was One) -> WazirX
$ws.Range("B9").Value = "WazirX"
$ws.Range("C9").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextCell $ws.Range("D9") "0.1428"
$ws.Range("E9").Value = "8WazirXWRX"

# Row 10 (was WazirX) -> MandalaExchangeToken
$ws.Range("B10").Value = "MandalaExchangeToken"
$ws.Range("C10").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextCell $ws.Range("D10") "0.07432"
$ws.Range("E10").Value = "9MandalaExchangeTokenMDX"

# Row 11 (was MandalaExchangeToken) -> LiechtensteinCryptoassetsExchange
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextCell $ws.Range("D11") "0.03378"
$ws.Range("E11").Value = "10LiechtensteinCryptoassetsExchangeLCX"

# Row 12 (was LiechtensteinCryptoassetsExchange) -> BitrueCoin
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextCell $ws.Range("D12") "0.03035"
$ws.Range("E12").Value = "11BitrueCoinBTR"

# Row 13 (was BitrueCoin) -> BitMartToken
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextCell $ws.Range("D13") "0.09358"
$ws.Range("E13").Value = "12BitMartTokenBMX"

# Row 14 (was BitMartToken) -> MCDex
$ws.Range("B14").Value = "MCDex"
$ws.Range("C14").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextCell $ws.Range("D14") "3.945"
$ws.Range("E14").Value = "13MCDexMCB"

# Row 15 (was MCDex) -> BitForexToken
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextCell $ws.Range("D15") "0.001596"
$ws.Range("E15").Value = "14BitForexTokenBF"

# Row 16 (was BitForexToken) -> CoinExToken
$ws.Range("B16").Value = "CoinExToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextCell $ws.Range("D16") "0.04805"
$ws.Range("E16").Value = "15CoinExTokenCET"

# Row 17 (was CoinExToken) -> One
$ws.Range("B17").Value = "One"
$ws.Range("C17").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextCell $ws.Range("D17") "0.0005942"
$ws.Range("E17").Value = "16OneONE"

# ---------------------------------------------------------------------
# More simple "Price" refreshes (column D only)
# ---------------------------------------------------------------------
Set-TextCell $ws.Range("D18") "0.005481"
Set-TextCell $ws.Range("D20") "0.0009870"
Set-TextCell $ws.Range("D21") "0.00007703"
Set-TextCell $ws.Range("D22") "3.660"
Set-TextCell $ws.Range("D23") "6.458"
Set-TextCell $ws.Range("D40") "0.03951"

# ---------------------------------------------------------------------
# Rows 41-43: KickToken / BKEXToken / CEJI rotate positions.
# ---------------------------------------------------------------------

# Row 41 (was KickToken) -> BKEXToken
$ws.Range("B41").Value = "BKEXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextCell $ws.Range("D41") "0.1075"
$ws.Range("E41").Value = "40BKEXTokenBKK"

# Row 42 (was BKEXToken) -> CEJI
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextCell $ws.Range("D42") "0.002711"
$ws.Range("E42").Value = "41CEJICEJI"

# Row 43 (was CEJI) -> KickToken
$ws.Range("B43").Value = "KickToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextCell $ws.Range("D43") "0.003029"
$ws.Range("E43").Value = "42KickTokenKICK"

# ---------------------------------------------------------------------
# Final simple "Price" refreshes + one Volume(1h) text tweak
# ---------------------------------------------------------------------
Set-TextCell $ws.Range("D44") "0.006926"
Set-TextCell $ws.Range("D45") "0.00005206"
$ws.Range("E48").Value = "47CoinbaseStockTokenCOINBestin24h"
